# Bond dates update: advance "today" by one day.
# For every data row (2..262):
#   - Column G ("Dni od poprzedniej wyplaty" / days since previous payment) increases by 1,
#     but only when the cell actually holds a value (some rows have no previous payment date).
#   - Column I ("Dni do nastepnej wyplaty" / days to next payment) decreases by 1 for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 262

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $gVal = $gCell.Value2
    if ($gVal -ne $null) {
        $gCell.Value2 = $gVal + 1
    }

    $iCell = $ws.Cells.Item($r, 9)   # column I
    $iVal = $iCell.Value2
    if ($iVal -ne $null) {
        $iCell.Value2 = $iVal - 1
    }
}
